# Raw and Clean Data from SSA for June 10th
# Adds a new daily-tracking row (row 10) to the out_vars sheet and
# normalizes the date column's formatting/width, mirroring what happened
# when the workbook was reopened/edited and resaved in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring down the same formatting (bold header-less body style + borders)
# that row 9 already has, so the new row looks consistent with the rest
# of the table, then overwrite with the new day's values.
$ws.Range("A9:G9").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 43991
$ws.Range("C10").Value = 124301
$ws.Range("D10").Value = 182077
$ws.Range("E10").Value = 50677
$ws.Range("F10").Value = 14649
$ws.Range("G10").Value = 33.21

# The date column previously mixed two slightly different custom date
# formats (with/without a time component). Normalize the whole column,
# including the freshly added row, onto a single short-date format.
$ws.Range("B2:B10").NumberFormat = "m/d/yy"

# Widen column B so the dates are fully visible (matches the workbook's
# best-fit column width after the edit).
$ws.Columns.Item(2).ColumnWidth = 16.8307

# Leave selection on the last cell that was touched, like a human editor
# would after typing in the new row.
$ws.Range("F10").Select() | Out-Null
